$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Elementos"
# ----------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elementos")

# Apply the bordered style (already used by A8, which has style s="1") to A2:A7
$fmtSrc = $ws.Range("A8")
$fmtSrc.Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)

# Remove the now-unused B/C columns for rows 8-22 (content was already empty,
# only formatting was present) - clear formats first, then contents, so the
# cell records disappear completely from the sheet.
$ws.Range("B8:C22").ClearFormats()
$ws.Range("B8:C22").ClearContents()

$ws.Range("E20").Select()
